$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = "Gilberto-Mec. Tec. Res. M"

$ws.Range("C11").Value = "Ludoff-Máquinas Térmicas e de Fl"
$ws.Range("D11").Value = "Gilberto-Mec. Tec. Res. M"
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "[Anderson-Ajustagem, Aline S. M.-Metalografia, Aderci-Tornearia, Andre B.-Elet. Dig. Bas.]"

$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "[Suzanny-Des. Maq. Cad_T1, Claudinei-Des. Maq. Cad_T2]"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "[Anderson-Ajustagem, Aline S. M.-Metalografia, Aderci-Tornearia, Andre B.-Elet. Dig. Bas.]"

$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "[Suzanny-Des. Maq. Cad_T1, Claudinei-Des. Maq. Cad_T2]"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "[Anderson-Ajustagem, Aline S. M.-Metalografia, Aderci-Tornearia, Andre B.-Elet. Dig. Bas.]"

$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "[Suzanny-Des. Maq. Cad_T1, Claudinei-Des. Maq. Cad_T2]"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "[Anderson-Ajustagem, Aline S. M.-Metalografia, Aderci-Tornearia, Andre B.-Elet. Dig. Bas.]"

$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "Ludoff-Máquinas Térmicas e de Fl"
